# Auto-generated Excel COM-interop script to update cryptos price/volume data
# (reflects scraped price/volume refresh from GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text such as '212.07' or '26.640.30'. Excel's COM
# Value setter auto-parses single-dot numeric-looking text into a Double,
# which would corrupt values (e.g. 212.07 -> 212.06999999999999) and change
# the cell type away from text. Force the Price column to Text format while
# we write the values, then restore the original (default/General) style so
# the cells keep matching their original 'no explicit style' formatting.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.640.30'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '1.595.20'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '212.07'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').Value = '0.512'
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -1.61%  '
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('D10').Value = '19.66'
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('D11').Value = '0.0836'
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('D12').Value = '1.819.71'
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('D13').Value = '1.595.88'
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  -2.78%  '
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  -2.88%  '
$ws.Range('D16').Value = '65.15'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('D17').Value = '26.623.03'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  -2.20%  '
$ws.Range('D19').Value = '209.79'
$ws.Range('E19').Value = '  -1.90%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = '6.70'
$ws.Range('E21').Value = '  -2.28%  '
$ws.Range('D22').Value = '4.25'
$ws.Range('E22').Value = '  -2.51%  '
$ws.Range('D23').Value = '2.32'
$ws.Range('E23').Value = '  -2.49%  '
$ws.Range('D24').Value = '8.88'
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('D25').Value = '146.36'
$ws.Range('E25').Value = '  -1.24%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').Value = '7.15'
$ws.Range('E27').Value = '  -3.03%  '
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('E29').Value = '  -1.34%  '
$ws.Range('D30').Value = '0.0505'
$ws.Range('E30').Value = '  -1.67%  '
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').Value = '3.23'
$ws.Range('E32').Value = '  -3.40%  '
$ws.Range('D33').Value = '0.667'
$ws.Range('E33').Value = '  -12.46%  '
$ws.Range('E34').Value = '  -3.29%  '
$ws.Range('D35').Value = '1.295.85'
$ws.Range('E35').Value = '  -3.86%  '
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('D37').Value = '1.49'
$ws.Range('E37').Value = '  -4.86%  '
$ws.Range('E38').Value = '  -3.16%  '
$ws.Range('D39').Value = '0.834'
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('E42').Value = '  +0.73%  '
$ws.Range('D43').Value = '2.20'
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('D44').Value = '63.47'
$ws.Range('E44').Value = '  -2.48%  '
$ws.Range('D45').Value = '1.731.90'
$ws.Range('E45').Value = '  -1.52%  '
$ws.Range('D46').Value = '89.42'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('D48').Value = '0.829'
$ws.Range('E48').Value = '  -5.68%  '
$ws.Range('E49').Value = '  -3.24%  '
$ws.Range('E50').Value = '  -2.45%  '
$ws.Range('D51').Value = '7.56'
$ws.Range('E51').Value = '  -1.10%  '

# Restore the Price column formatting/style so no stray explicit style index
# is left behind on cells that originally had none.
$priceRange.NumberFormat = "General"
$priceRange.Style = "Normal"
